# Update countries & provincias Spain
#
# 1) "Tayikistan" moves in the country list: it used to be the last entry
#    before "Guam" (row 180), and becomes the first entry right after
#    "Mozambique" (row 164), pushing the rows that used to be 164-180
#    down to 165-181... but since Tayikistan is removed from its old
#    spot, everything nets out and rows 165-180 simply take on the
#    country name + stats that used to belong to the row above them
#    (164..179), while row 164 gets brand-new data for Tayikistan.
# 2) A handful of unrelated rows get refreshed case/province counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Part 1: simple value refreshes (no row movement)
# ---------------------------------------------------------------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1149878
$ws.Range("C4").Value = 18848
$ws.Range("D4").Value = 162862
$ws.Range("E4").Value = 920197
$ws.Range("G4").Value = 1066
$ws.Range("H4").Value = 66819

# Row 9: Alemania
$ws.Range("B9").Value = 164602
$ws.Range("C9").Value = 525
$ws.Range("E9").Value = 28847
$ws.Range("G9").Value = 19
$ws.Range("H9").Value = 6755

# Row 19: India
$ws.Range("B19").Value = 39699
$ws.Range("C19").Value = 2442
$ws.Range("D19").Value = 10819
$ws.Range("E19").Value = 27557
$ws.Range("G19").Value = 100
$ws.Range("H19").Value = 1323

# Row 20: Suiza
$ws.Range("E20").Value = 4155
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 1762

# Row 133: Ruanda
$ws.Range("B133").Value = 255
$ws.Range("C133").Value = 6
$ws.Range("D133").Value = 120
$ws.Range("E133").Value = 135

# ---------------------------------------------------------------------
# Part 2: move "Tayikistan" from its old row (180) to right after
# "Mozambique" (new row 164), shifting the rows in between down by one.
# ---------------------------------------------------------------------

# Capture the old country names + stats for rows 164..179 (these are the
# rows that must shift down to 165..180). Use .Value2 for reads - in this
# runtime, a bare ".Value" getter returns the property descriptor instead
# of invoking it, while ".Value2" reads the actual cell content.
$oldRows = @()
for ($r = 164; $r -le 179; $r++) {
    $oldRows += ,@(
        $ws.Cells.Item($r, 1).Value2,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2,
        $ws.Cells.Item($r, 7).Value2,
        $ws.Cells.Item($r, 8).Value2
    )
}

# Write them back starting at row 165 (shift down by one row).
for ($i = 0; $i -lt $oldRows.Count; $i++) {
    $r = 165 + $i
    $row = $oldRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Now write the new Tayikistan row at 164 with its refreshed stats.
$ws.Cells.Item(164, 1).Value = "Tayikistan"
$ws.Cells.Item(164, 2).Value = 76
$ws.Cells.Item(164, 3).Value = 61
$ws.Cells.Item(164, 4).Value = 0
$ws.Cells.Item(164, 5).Value = 74
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 2
$ws.Cells.Item(164, 8).Value = 2
